$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting the existing "Tipo"/"single" column to E
$ws.Columns.Item(4).Insert()

# New header cell D1 = "MAE" (inherits header formatting from the column insert)
$ws.Range("D1").Value = "MAE"

# New data cell D2 = MAE value
$ws.Range("D2").Value = 0.1292627146720763
